# Auto-generated Excel COM-interop script
# Applies updated market board price data to the Leve profit sheets
# (commit: chore: update Sheets via scheduled runner)

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 19: Unbreak My Heart
$ws.Cells.Item(19, 8).Value = 431.95834
$ws.Cells.Item(19, 9).Value = 283.33334
$ws.Cells.Item(19, 10).Value = 521.13336
$ws.Cells.Item(19, 11).Value = 283.33334
$ws.Cells.Item(19, 12).Value = 521.13336
$ws.Cells.Item(19, 13).Value = -108.33334
$ws.Cells.Item(19, 14).Value = -871.13336

# Row 33: Glazed and Confused
$ws.Cells.Item(33, 8).Value = 321.77777
$ws.Cells.Item(33, 9).Value = 280.14285
$ws.Cells.Item(33, 11).Value = 280.14285
$ws.Cells.Item(33, 13).Value = -51.14285000000001

# Row 70: Consecrating Congregation
$ws.Cells.Item(70, 8).Value = 903.1539
$ws.Cells.Item(70, 10).Value = 1054.875
$ws.Cells.Item(70, 12).Value = 3164.625
$ws.Cells.Item(70, 14).Value = -3704.625

# Row 73: Curbing the Contagion (L)
$ws.Cells.Item(73, 8).Value = 903.1539
$ws.Cells.Item(73, 10).Value = 1054.875
$ws.Cells.Item(73, 12).Value = 3164.625
$ws.Cells.Item(73, 14).Value = -5036.625

# Row 74: Adhesive of Antipathy
$ws.Cells.Item(74, 8).Value = 6755.273
$ws.Cells.Item(74, 9).Value = 35000
$ws.Cells.Item(74, 11).Value = 35000
$ws.Cells.Item(74, 13).Value = -34064

# Row 76: Warding Off Temptation
$ws.Cells.Item(76, 8).Value = 3284.0476
$ws.Cells.Item(76, 9).Value = 3258.7646
$ws.Cells.Item(76, 11).Value = 3258.7646
$ws.Cells.Item(76, 13).Value = -2943.7646

# Row 77: It's Gonna Grow Back (L)
$ws.Cells.Item(77, 8).Value = 6755.273
$ws.Cells.Item(77, 9).Value = 35000
$ws.Cells.Item(77, 11).Value = 175000
$ws.Cells.Item(77, 13).Value = -170320

# Row 79: The Garden of Arcane Delights (L)
$ws.Cells.Item(79, 8).Value = 3284.0476
$ws.Cells.Item(79, 9).Value = 3258.7646
$ws.Cells.Item(79, 11).Value = 3258.7646
$ws.Cells.Item(79, 13).Value = -2166.7646

# Row 86: Filling in the Blanks
$ws.Cells.Item(86, 8).Value = 2346
$ws.Cells.Item(86, 9).Value = 1558.8
$ws.Cells.Item(86, 11).Value = 1558.8
$ws.Cells.Item(86, 13).Value = -435.8

# Row 87: There Was a Late Fee
$ws.Cells.Item(87, 8).Value = 33236
$ws.Cells.Item(87, 10).Value = 33236
$ws.Cells.Item(87, 12).Value = 33236
$ws.Cells.Item(87, 14).Value = -35732

# Row 89: Ink into Antiquity (L)
$ws.Cells.Item(89, 8).Value = 2346
$ws.Cells.Item(89, 9).Value = 1558.8
$ws.Cells.Item(89, 11).Value = 7794
$ws.Cells.Item(89, 13).Value = -2178

# Row 90: A Gate Arcane Is Dragon's Bane (L)
$ws.Cells.Item(90, 8).Value = 33236
$ws.Cells.Item(90, 10).Value = 33236
$ws.Cells.Item(90, 12).Value = 99708
$ws.Cells.Item(90, 14).Value = -112188

# Row 129: Practical Command
$ws.Cells.Item(129, 8).Value = 758.09375
$ws.Cells.Item(129, 9).Value = 521.4375
$ws.Cells.Item(129, 10).Value = 994.75
$ws.Cells.Item(129, 11).Value = 1564.3125
$ws.Cells.Item(129, 12).Value = 2984.25
$ws.Cells.Item(129, 13).Value = 3435.6875
$ws.Cells.Item(129, 14).Value = -12984.25

# Row 141: Remedy for Reason
$ws.Cells.Item(141, 8).Value = 3321.3333
$ws.Cells.Item(141, 9).Value = 2939.818
$ws.Cells.Item(141, 11).Value = 8819.454000000002
$ws.Cells.Item(141, 13).Value = -3639.454000000002

$ws = $wb.Worksheets.Item("ARM")
# Row 32: Ingot We Trust
$ws.Cells.Item(32, 8).Value = 9541.788
$ws.Cells.Item(32, 9).Value = 5924.924
$ws.Cells.Item(32, 10).Value = 23828.4
$ws.Cells.Item(32, 11).Value = 5924.924
$ws.Cells.Item(32, 12).Value = 23828.4
$ws.Cells.Item(32, 13).Value = -5637.924
$ws.Cells.Item(32, 14).Value = -24402.4

# Row 119: Trial and Error
$ws.Cells.Item(119, 8).Value = 29465.334
$ws.Cells.Item(119, 10).Value = 29465.334
$ws.Cells.Item(119, 12).Value = 29465.334
$ws.Cells.Item(119, 14).Value = -39141.334

$ws = $wb.Worksheets.Item("BSM")
# Row 22: Riveting Run
$ws.Cells.Item(22, 8).Value = 570.2857
$ws.Cells.Item(22, 9).Value = 570.2857
$ws.Cells.Item(22, 11).Value = 570.2857
$ws.Cells.Item(22, 13).Value = -397.2857

# Row 86: Through Thick and Thin
$ws.Cells.Item(86, 8).Value = 11593.125
$ws.Cells.Item(86, 9).Value = 14224.5
$ws.Cells.Item(86, 10).Value = 3699
$ws.Cells.Item(86, 11).Value = 14224.5
$ws.Cells.Item(86, 12).Value = 3699
$ws.Cells.Item(86, 13).Value = -13101.5
$ws.Cells.Item(86, 14).Value = -5945

# Row 89: Piercing Eyes Deserve Piercing Shafts (L)
$ws.Cells.Item(89, 8).Value = 11593.125
$ws.Cells.Item(89, 9).Value = 14224.5
$ws.Cells.Item(89, 10).Value = 3699
$ws.Cells.Item(89, 11).Value = 71122.5
$ws.Cells.Item(89, 12).Value = 18495
$ws.Cells.Item(89, 13).Value = -65506.5
$ws.Cells.Item(89, 14).Value = -29727

$ws = $wb.Worksheets.Item("CRP")
# Row 31: Wall Not Found
$ws.Cells.Item(31, 8).Value = 19365.797
$ws.Cells.Item(31, 9).Value = 1714.7742
$ws.Cells.Item(31, 10).Value = 33765.316
$ws.Cells.Item(31, 11).Value = 1714.7742
$ws.Cells.Item(31, 12).Value = 33765.316
$ws.Cells.Item(31, 13).Value = -1419.7742
$ws.Cells.Item(31, 14).Value = -34355.316

# Row 34: Armoires of the Rich and Famous
$ws.Cells.Item(34, 8).Value = 19365.797
$ws.Cells.Item(34, 9).Value = 1714.7742
$ws.Cells.Item(34, 10).Value = 33765.316
$ws.Cells.Item(34, 11).Value = 1714.7742
$ws.Cells.Item(34, 12).Value = 33765.316
$ws.Cells.Item(34, 13).Value = -1512.7742
$ws.Cells.Item(34, 14).Value = -34169.316

# Row 99: O Pine
$ws.Cells.Item(99, 8).Value = 1373.6
$ws.Cells.Item(99, 9).Value = 1267
$ws.Cells.Item(99, 10).Value = 1800
$ws.Cells.Item(99, 11).Value = 1267
$ws.Cells.Item(99, 12).Value = 1800
$ws.Cells.Item(99, 13).Value = 231
$ws.Cells.Item(99, 14).Value = -4796

# Row 126: A Better Conductor
$ws.Cells.Item(126, 8).Value = 1373.6
$ws.Cells.Item(126, 9).Value = 1267
$ws.Cells.Item(126, 10).Value = 1800
$ws.Cells.Item(126, 11).Value = 3801
$ws.Cells.Item(126, 12).Value = 5400
$ws.Cells.Item(126, 13).Value = -1331
$ws.Cells.Item(126, 14).Value = -10340

$ws = $wb.Worksheets.Item("CUL")
# Row 12: Butter Me Up
$ws.Cells.Item(12, 8).Value = 118.52941
$ws.Cells.Item(12, 9).Value = 74.066666
$ws.Cells.Item(12, 10).Value = 153.63158
$ws.Cells.Item(12, 11).Value = 222.199998
$ws.Cells.Item(12, 12).Value = 460.8947400000001
$ws.Cells.Item(12, 13).Value = -49.19999799999999
$ws.Cells.Item(12, 14).Value = -806.8947400000001

# Row 113: Can't Eat Just One
$ws.Cells.Item(113, 8).Value = 556.4737
$ws.Cells.Item(113, 9).Value = 510.0357
$ws.Cells.Item(113, 10).Value = 601.3103599999999
$ws.Cells.Item(113, 11).Value = 1530.1071
$ws.Cells.Item(113, 12).Value = 1803.93108
$ws.Cells.Item(113, 13).Value = 639.8928999999998
$ws.Cells.Item(113, 14).Value = -6143.93108

# Row 131: The Mountain Steeped
$ws.Cells.Item(131, 8).Value = 905.4706
$ws.Cells.Item(131, 10).Value = 996.11365
$ws.Cells.Item(131, 12).Value = 2988.34095
$ws.Cells.Item(131, 14).Value = -13068.34095

$ws = $wb.Worksheets.Item("GSM")
# Row 80: Needs More Prayerbell
$ws.Cells.Item(80, 8).Value = 4571.4165
$ws.Cells.Item(80, 9).Value = 4365.3335
$ws.Cells.Item(80, 10).Value = 4640.1113
$ws.Cells.Item(80, 11).Value = 4365.3335
$ws.Cells.Item(80, 12).Value = 4640.1113
$ws.Cells.Item(80, 13).Value = -3367.3335
$ws.Cells.Item(80, 14).Value = -6636.1113

# Row 83: With a Noise That Reaches Heaven (L)
$ws.Cells.Item(83, 8).Value = 4571.4165
$ws.Cells.Item(83, 9).Value = 4365.3335
$ws.Cells.Item(83, 10).Value = 4640.1113
$ws.Cells.Item(83, 11).Value = 21826.6675
$ws.Cells.Item(83, 12).Value = 23200.5565
$ws.Cells.Item(83, 13).Value = -16834.6675
$ws.Cells.Item(83, 14).Value = -33184.5565

$ws = $wb.Worksheets.Item("LTW")
# Row 82: Trainin' the Neck
$ws.Cells.Item(82, 8).Value = 1217.5555
$ws.Cells.Item(82, 9).Value = 1082.25
$ws.Cells.Item(82, 10).Value = 2300
$ws.Cells.Item(82, 11).Value = 1082.25
$ws.Cells.Item(82, 12).Value = 2300
$ws.Cells.Item(82, 13).Value = -721.25
$ws.Cells.Item(82, 14).Value = -3022

# Row 85: Training Is Only Skintight (L)
$ws.Cells.Item(85, 8).Value = 1217.5555
$ws.Cells.Item(85, 9).Value = 1082.25
$ws.Cells.Item(85, 10).Value = 2300
$ws.Cells.Item(85, 11).Value = 1082.25
$ws.Cells.Item(85, 12).Value = 2300
$ws.Cells.Item(85, 13).Value = 165.75
$ws.Cells.Item(85, 14).Value = -4796

# Row 93: Hide to Go Seek
$ws.Cells.Item(93, 8).Value = 1111.6666
$ws.Cells.Item(93, 9).Value = 1128.4
$ws.Cells.Item(93, 10).Value = 777
$ws.Cells.Item(93, 11).Value = 1128.4
$ws.Cells.Item(93, 12).Value = 777
$ws.Cells.Item(93, 13).Value = 119.5999999999999
$ws.Cells.Item(93, 14).Value = -3273

# Row 98: Try Tricorne Again
$ws.Cells.Item(98, 8).Value = 28427
$ws.Cells.Item(98, 10).Value = 28427
$ws.Cells.Item(98, 12).Value = 28427
$ws.Cells.Item(98, 14).Value = -34417

# Row 100: Tiger in the Sack
$ws.Cells.Item(100, 8).Value = 1873.4
$ws.Cells.Item(100, 9).Value = 1765.125
$ws.Cells.Item(100, 10).Value = 1997.1428
$ws.Cells.Item(100, 11).Value = 1765.125
$ws.Cells.Item(100, 12).Value = 1997.1428
$ws.Cells.Item(100, 13).Value = -1224.125
$ws.Cells.Item(100, 14).Value = -3079.1428

# Row 136: Respect for Br'aax
$ws.Cells.Item(136, 8).Value = 64430.062
$ws.Cells.Item(136, 9).Value = 42281.04
$ws.Cells.Item(136, 10).Value = 146697.86
$ws.Cells.Item(136, 11).Value = 126843.12
$ws.Cells.Item(136, 12).Value = 440093.58
$ws.Cells.Item(136, 13).Value = -124293.12
$ws.Cells.Item(136, 14).Value = -445193.58

$ws = $wb.Worksheets.Item("WVR")
# Row 107: Flax Wax
$ws.Cells.Item(107, 8).Value = 395.26666
$ws.Cells.Item(107, 9).Value = 319.96875
$ws.Cells.Item(107, 10).Value = 580.61536
$ws.Cells.Item(107, 11).Value = 959.90625
$ws.Cells.Item(107, 12).Value = 1741.84608
$ws.Cells.Item(107, 13).Value = 960.09375
$ws.Cells.Item(107, 14).Value = -5581.84608

# Row 112: Hair Do No Harm
$ws.Cells.Item(112, 8).Value = 39741
$ws.Cells.Item(112, 10).Value = 39741
$ws.Cells.Item(112, 12).Value = 39741
$ws.Cells.Item(112, 14).Value = -42695

# Row 113: A Tender Table
$ws.Cells.Item(113, 8).Value = 629.80554
$ws.Cells.Item(113, 9).Value = 776.0476
$ws.Cells.Item(113, 10).Value = 425.06668
$ws.Cells.Item(113, 11).Value = 2328.1428
$ws.Cells.Item(113, 12).Value = 1275.20004
$ws.Cells.Item(113, 13).Value = -158.1428000000001
$ws.Cells.Item(113, 14).Value = -5615.20004

# Row 132: Comfy Cabins
$ws.Cells.Item(132, 8).Value = 75354.81
$ws.Cells.Item(132, 9).Value = 59637.176
$ws.Cells.Item(132, 10).Value = 102074.8
$ws.Cells.Item(132, 11).Value = 178911.528
$ws.Cells.Item(132, 12).Value = 306224.4
$ws.Cells.Item(132, 13).Value = -176381.528
$ws.Cells.Item(132, 14).Value = -311284.4
